$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'59.676.33"
$ws.Cells.Item(2, 5).Value = "  +2.31%  "
$ws.Cells.Item(3, 4).Value = "'2.416.97"
$ws.Cells.Item(3, 5).Value = "  +2.54%  "
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$ws.Cells.Item(5, 4).Value = "'554.27"
$ws.Cells.Item(5, 5).Value = "  +2.53%  "
$ws.Cells.Item(6, 4).Value = "'138.10"
$ws.Cells.Item(6, 5).Value = "  +1.88%  "
$ws.Cells.Item(7, 5).Value = "  -0.21%  "
$ws.Cells.Item(8, 5).Value = "  +0.93%  "
$ws.Cells.Item(9, 4).Value = "'0.107"
$ws.Cells.Item(9, 5).Value = "  +4.21%  "
$ws.Cells.Item(10, 4).Value = "'5.83"
$ws.Cells.Item(10, 5).Value = "  +4.41%  "
$ws.Cells.Item(11, 5).Value = "  +1.52%  "
$ws.Cells.Item(12, 5).Value = "  -2.14%  "
$ws.Cells.Item(13, 4).Value = "'24.71"
$ws.Cells.Item(13, 5).Value = "  +3.74%  "
$ws.Cells.Item(14, 4).Value = "'2.850.02"
$ws.Cells.Item(14, 5).Value = "  +2.67%  "
$ws.Cells.Item(15, 4).Value = "'59.608.44"
$ws.Cells.Item(15, 5).Value = "  +2.20%  "
$ws.Cells.Item(16, 5).Value = "  +4.06%  "
$ws.Cells.Item(17, 4).Value = "'2.435.93"
$ws.Cells.Item(17, 5).Value = "  +3.30%  "
$ws.Cells.Item(18, 4).Value = "'11.37"
$ws.Cells.Item(18, 5).Value = "  +5.99%  "
$ws.Cells.Item(20, 4).Value = "'333.23"
$ws.Cells.Item(20, 5).Value = "  +0.07%  "
$ws.Cells.Item(21, 4).Value = "'6.91"
$ws.Cells.Item(21, 5).Value = "  +1.71%  "
$ws.Cells.Item(22, 5).Value = "  -0.06%  "
$ws.Cells.Item(23, 4).Value = "'64.51"
$ws.Cells.Item(23, 5).Value = "  +2.32%  "
$ws.Cells.Item(24, 5).Value = "  +1.41%  "
$ws.Cells.Item(25, 4).Value = "'8.57"
$ws.Cells.Item(25, 5).Value = "  +1.24%  "
$ws.Cells.Item(26, 5).Value = "  +0.11%  "
$ws.Cells.Item(27, 5).Value = "  -1.72%  "
$ws.Cells.Item(28, 4).Value = "'0.0₃0784"
$ws.Cells.Item(28, 5).Value = "  +6.69%  "
$ws.Cells.Item(29, 5).Value = "  +3.41%  "
$ws.Cells.Item(30, 4).Value = "'170.71"
$ws.Cells.Item(30, 5).Value = "  -0.74%  "
$ws.Cells.Item(31, 5).Value = "  +1.97%  "
$ws.Cells.Item(32, 4).Value = "'18.64"
$ws.Cells.Item(32, 5).Value = "  +1.02%  "
$ws.Cells.Item(33, 5).Value = "  -1.44%  "
$ws.Cells.Item(34, 5).Value = "  -0.02%  "
$ws.Cells.Item(35, 4).Value = "'1.31"
$ws.Cells.Item(35, 5).Value = "  +5.28%  "
$ws.Cells.Item(36, 4).Value = "'4.24"
$ws.Cells.Item(36, 5).Value = "  -0.93%  "
$ws.Cells.Item(37, 5).Value = "  +0.15%  "
$ws.Cells.Item(38, 5).Value = "  -0.85%  "
$ws.Cells.Item(39, 4).Value = "'40.10"
$ws.Cells.Item(39, 5).Value = "  +2.18%  "
$ws.Cells.Item(40, 4).Value = "'0.422"
$ws.Cells.Item(40, 5).Value = "  +11.50%  "
$ws.Cells.Item(41, 4).Value = "'312.99"
$ws.Cells.Item(41, 5).Value = "  +6.09%  "
$ws.Cells.Item(42, 5).Value = "  +2.13%  "
$ws.Cells.Item(43, 4).Value = "'142.93"
$ws.Cells.Item(43, 5).Value = "  -1.93%  "
$ws.Cells.Item(44, 4).Value = "'0.0962"
$ws.Cells.Item(44, 5).Value = "  +1.52%  "
$ws.Cells.Item(45, 4).Value = "'0.0524"
$ws.Cells.Item(45, 5).Value = "  +4.16%  "
$ws.Cells.Item(46, 4).Value = "'0.409"
$ws.Cells.Item(46, 5).Value = "  +6.05%  "
$ws.Cells.Item(47, 4).Value = "'19.15"
$ws.Cells.Item(47, 5).Value = "  -0.22%  "
$ws.Cells.Item(48, 5).Value = "  +1.24%  "
$ws.Cells.Item(49, 5).Value = "  +2.49%  "
$ws.Cells.Item(50, 4).Value = "'11.05"
$ws.Cells.Item(50, 5).Value = "  -0.28%  "
$ws.Cells.Item(51, 5).Value = "  +4.15%  "
